$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lapa1")

# Delete rows 4 and 5, which shifts the old rows 6 and 7 up into their place.
$ws.Rows("4:5").Delete()
